# approvalProcess.xlsx - "Updated the test data."
# Adds a new "supplier"/"Supplier" column (V) to the FinanceApprover sheet,
# (re)sets the best-fit column widths for the header row, and updates the
# current selection to the full S column (clearing the old scrolled view).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data column: V1 = "supplier" header, V2 = "Supplier" value ---
$ws.Range("V1").Value = "supplier"
$ws.Range("V2").Value = "Supplier"

# --- Column widths (best-fit sizing after the new column was added) ---
$ws.Columns("A").ColumnWidth  = 10.592447916666666
$ws.Columns("B").ColumnWidth  = 9.451822916666666
$ws.Columns("C").ColumnWidth  = 16.451822916666668
$ws.Columns("D").ColumnWidth  = 12.307291666666666
$ws.Columns("E").ColumnWidth  = 20.307291666666668
$ws.Columns("F").ColumnWidth  = 26.022135416666668
$ws.Columns("G").ColumnWidth  = 10.451822916666666
$ws.Columns("H").ColumnWidth  = 11.736979166666666
$ws.Columns("I").ColumnWidth  = 18.166666666666668
$ws.Columns("J").ColumnWidth  = 24.307291666666668
$ws.Columns("K").ColumnWidth  = 7.877604166666667
$ws.Columns("L").ColumnWidth  = 17.166666666666668
$ws.Columns("M").ColumnWidth  = 5.022135416666667
$ws.Columns("P").ColumnWidth  = 5.736979166666667
$ws.Columns("R").ColumnWidth  = 14.307291666666666
$ws.Columns("S").ColumnWidth  = 9.592447916666666
$ws.Columns("T").ColumnWidth  = 28.022135416666668
$ws.Columns("U").ColumnWidth  = 13.877604166666666
$ws.Columns("V").ColumnWidth  = 7.592447916666667

# --- Selection: whole column S is selected, and the view is scrolled back
#     so there is no stored topLeftCell anymore ---
$ws.Columns("S").Select()
